$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Battle")

# Insert two new rows above the current row 12 ("Question:"), shifting the
# rest of the battle-card layout (question box, DATA/X labels, blank grid)
# down by two rows.
$ws.Rows("12:13").Insert()

# Fill the newly inserted row 12 with the "give up" instruction.
$ws.Range("B12").Value = "If you give up, write '-1' as the answer."
